# Fix Haggai reference cells in column A: strip the erroneous trailing "16"
# that was accidentally appended to many of the verse references, so they
# read as clean "Book Chapter:Verse" strings for human readability.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixes = @{
    "A2"  = "Haggai 1:2"
    "A3"  = "Haggai 1:3"
    "A4"  = "Haggai 1:5"
    "A5"  = "Haggai 1:7"
    "A7"  = "Haggai 1:10"
    "A8"  = "Haggai 1:11"
    "A9"  = "Haggai 1:13"
    "A10" = "Haggai 1:14"
    "A11" = "Haggai 2:2"
    "A12" = "Haggai 2:4"
    "A14" = "Haggai 2:7"
    "A15" = "Haggai 2:11"
    "A16" = "Haggai 2:12"
    "A17" = "Haggai 2:13"
    "A18" = "Haggai 2:14"
    "A19" = "Haggai 2:15"
    "A20" = "Haggai 2:16"
    "A21" = "Haggai 2:20"
    "A22" = "Haggai 2:21"
    "A23" = "Haggai 2:22"
    "A24" = "Haggai 2:23"
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}
